$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (EXP001) updates
$ws.Range("B2").Value = "n1000000_f_init5_cont0_disc5_sep5p1_seed0"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "2025-07-27"
$ws.Range("G2").Value = 0.99995754
$ws.Range("I2").Value = "reports/n1000000_f_init5_cont0_disc5_sep5p1_seed0_mlp_001_tuning.db"
$ws.Range("J2").Value = "reports/figures/n1000000_f_init5_cont0_disc5_sep5p1_seed0/mlp"

# Row 3 (EXP002) updates
$ws.Range("B3").Value = "n1000000_f_init5_cont0_disc5_sep5p1_seed0"
$ws.Range("D3").Value = "Final Training"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "2025-07-27"
$ws.Range("G3").Value = 0.99624
$ws.Range("H3").Value = "N/A"
$ws.Range("I3").Value = "models/n1000000_f_init5_cont0_disc5_sep5p1_seed0_mlp_001_optimal_model_metrics.json"
$ws.Range("J3").Value = "reports/figures/n1000000_f_init5_cont0_disc5_sep5p1_seed0/mlp_final"
